# Insert a new price-report row into the "Berenjena" sheet at row 508.
# This shifts all existing rows 508..542 down by one (to 509..543),
# and populates the newly opened row 508 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 508, pushing the rest down.
$ws.Rows.Item(508).Insert()

# Populate the new row 508 with the new data record.
$ws.Range("A508").Value = 10
$ws.Range("B508").Value = "Vega Modelo de Temuco"
$ws.Range("C508").Value = "La Araucanía"
$ws.Range("D508").Value = 45265
$ws.Range("E508").Value = 9
$ws.Range("F508").Value = 100112001
$ws.Range("G508").Value = "Berenjena"
$ws.Range("H508").Value = "Sin especificar"
$ws.Range("I508").Value = "Primera"
$ws.Range("J508").Value = 100
$ws.Range("K508").Value = 15000
$ws.Range("L508").Value = 15000
$ws.Range("M508").Value = 15000
$ws.Range("N508").Value = "$/caja 40 unidades"
$ws.Range("O508").Value = "Región de Arica y Parinacota"
$ws.Range("P508").Value = 375
$ws.Range("Q508").Value = 40
$ws.Range("R508").Value = "Hortaliza"
